$wb = $excel.ActiveWorkbook

# --- Sheet1 (Login) : add a "pageheader" column (C) with "Swag Labs" values ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C1").Value = "pageheader"
$ws1.Range("C2").Value = "Swag Labs"
$ws1.Range("C3").Value = "Swag Labs"
$ws1.Range("C4").Value = "Swag Labs"
$ws1.Range("C5").Value = "Swag Labs"
$ws1.Range("C6").Value = "Swag Labs"
$ws1.Range("C7").Value = "Swag Labs"

# Update the selection on the Login sheet
$null = $ws1.Range("A1:C2").Select()

# --- Sheet2 : rename to FilterProduct and add data rows ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "FilterProduct"

$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("C1").Value = "pageheader"
$ws2.Range("D1").Value = "filterproduct"

$ws2.Range("A2").Value = "standard_user"
$ws2.Range("B2").Value = "secret_sauce"
$ws2.Range("C2").Value = "Swag Labs"
$ws2.Range("D2").Value = "Price (low to high)"

# Set the selection on the FilterProduct sheet and make it the active tab
$null = $ws2.Range("D3").Select()
$ws2.Activate()
